$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44326, 1, 10, 107.7121930202499),
    @(44327, 1, 9, 96.9409737182249),
    @(44328, 0, 9, 96.9409737182249),
    @(44329, 0, 5, 53.85609651012494)
)

$lastExistingRow = 251
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy formatting of column A from the row immediately above so the
    # new date cell keeps the same style (border/font/number format).
    $srcA = $ws.Cells.Item($row - 1, 1)
    $dstA = $ws.Cells.Item($row, 1)
    $srcA.Copy()
    $dstA.PasteSpecial(-4122)  # xlPasteFormats

    $dstA.Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = $false
